# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-25, replacing the old "Strike#" derived values
$kValues = @{
    2  = 3
    3  = 2
    4  = 8
    5  = 6
    6  = 4
    7  = 9
    8  = 3
    9  = 6
    10 = 1
    11 = 4
    12 = 6
    13 = 6
    14 = 5
    15 = 8
    16 = 4
    17 = 9
    18 = 7
    19 = 4
    20 = 7
    21 = 1
    22 = 2
    23 = 6
    24 = 4
    25 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
